# fixed std error on before 26th map
# Row 25 computes the standard error (STDEV / sqrt(n)) of each column's
# stats in row 24. The sample-size divisor was wrong (71 instead of the
# actual number of observations, 17) - fix every cell in row 25 that
# divides by the old (71 ^ 0.5). (D25/G25 key off blank row-24 cells, so
# their result is 0 either way and they were left as-is.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data.csv")

$cols = @("B","C","E","F","H","I","L","M","N","O")
foreach ($col in $cols) {
    $cell = $col + "25"
    $src  = $col + "24"
    $ws.Range($cell).Formula = "=" + $src + "/(17 ^ 0.5)"
}

# Restore the view state recorded in the edited workbook: the active
# window was scrolled so row 11 is at the top, and the selected cell
# moved to O26.
$ws.Range("O26").Select() | Out-Null

$wb.Save() | Out-Null
